$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (2026-02-07 -> 2026-02-08)
# for every data row (rows 2 through 109). Update them all.
$range = $ws.Range("C2:C109")
$range.Value = 46061
